# Auto-generated cell updates applying the scheduled market-data refresh
# to the Chocobo_Profits leve-profit workbook (per-class sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 5
$ws.Range("H5").Value = 208.8
$ws.Range("I5").Value = 84
$ws.Range("K5").Value = 84
$ws.Range("M5").Value = 31
# row 41
$ws.Range("H41").Value = 446.25
$ws.Range("I41").Value = 130
$ws.Range("J41").Value = 636
$ws.Range("K41").Value = 130
$ws.Range("L41").Value = 636
$ws.Range("M41").Value = 310
$ws.Range("N41").Value = -1516
# row 92
$ws.Range("H92").Value = 1518.375
$ws.Range("I92").Value = 1315.75
$ws.Range("J92").Value = 2126.25
$ws.Range("K92").Value = 1315.75
$ws.Range("L92").Value = 2126.25
$ws.Range("M92").Value = -67.75
$ws.Range("N92").Value = -4622.25
# row 93
$ws.Range("H93").Value = 29624.172
$ws.Range("J93").Value = 29624.172
$ws.Range("L93").Value = 29624.172
$ws.Range("N93").Value = -34616.172
# row 137
$ws.Range("H137").Value = 2978911
$ws.Range("I137").Value = 3969923
$ws.Range("K137").Value = 11909769
$ws.Range("M137").Value = -11907219

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 102
$ws.Range("H102").Value = 2909.8
$ws.Range("I102").Value = 2642.7144
$ws.Range("J102").Value = 3533
$ws.Range("K102").Value = 2642.7144
$ws.Range("L102").Value = 3533
$ws.Range("M102").Value = -1020.7144
$ws.Range("N102").Value = -6777
# row 132
$ws.Range("H132").Value = 3170.4285
$ws.Range("I132").Value = 1259
$ws.Range("J132").Value = 3935
$ws.Range("K132").Value = 3777
$ws.Range("L132").Value = 11805
$ws.Range("M132").Value = -1247
$ws.Range("N132").Value = -16865
# row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 95
$ws.Range("H95").Value = 32663.158
$ws.Range("J95").Value = 32663.158
$ws.Range("L95").Value = 32663.158
$ws.Range("N95").Value = -38155.158
# row 99
$ws.Range("H99").Value = 4691
$ws.Range("I99").Value = 1495
$ws.Range("J99").Value = 5401.222
$ws.Range("K99").Value = 1495
$ws.Range("L99").Value = 5401.222
$ws.Range("M99").Value = 3
$ws.Range("N99").Value = -8397.222
# row 134
$ws.Range("H134").Value = 3030.9167
$ws.Range("I134").Value = 2313.7896
$ws.Range("K134").Value = 6941.3688
$ws.Range("M134").Value = -4406.3688

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2790.4243
$ws.Range("I31").Value = 915.05554
$ws.Range("J31").Value = 5040.8667
$ws.Range("K31").Value = 915.05554
$ws.Range("L31").Value = 5040.8667
$ws.Range("M31").Value = -620.05554
$ws.Range("N31").Value = -5630.8667
# row 34
$ws.Range("H34").Value = 2790.4243
$ws.Range("I34").Value = 915.05554
$ws.Range("J34").Value = 5040.8667
$ws.Range("K34").Value = 915.05554
$ws.Range("L34").Value = 5040.8667
$ws.Range("M34").Value = -713.05554
$ws.Range("N34").Value = -5444.8667
# row 58
$ws.Range("H58").Value = 3027.7424
$ws.Range("I58").Value = 1756.1296
$ws.Range("K58").Value = 1756.1296
$ws.Range("M58").Value = -1553.1296
# row 109
$ws.Range("H109").Value = 34999.668
$ws.Range("J109").Value = 34999.668
$ws.Range("L109").Value = 34999.668
$ws.Range("N109").Value = -37079.668
# row 134
$ws.Range("H134").Value = 2607.3572
$ws.Range("I134").Value = 1240
$ws.Range("J134").Value = 3367
$ws.Range("K134").Value = 3720
$ws.Range("L134").Value = 10101
$ws.Range("M134").Value = -1185
$ws.Range("N134").Value = -15171
# row 136
$ws.Range("H136").Value = 3027.7424
$ws.Range("I136").Value = 1756.1296
$ws.Range("K136").Value = 5268.3888
$ws.Range("M136").Value = -2718.3888

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 1269433.8
$ws.Range("I4").Value = 3444277.8
$ws.Range("J4").Value = 774.75
$ws.Range("K4").Value = 10332833.4
$ws.Range("L4").Value = 2324.25
$ws.Range("M4").Value = -10332721.4
$ws.Range("N4").Value = -2548.25
# row 22
$ws.Range("H22").Value = 2005.7333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2005.7333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6017.199900000001
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -6355.199900000001
# row 27
$ws.Range("H27").Value = 2005.7333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2005.7333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6017.199900000001
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = -6221.199900000001
# row 63
$ws.Range("H63").Value = 3950.7144
$ws.Range("I63").Value = 3434.4443
$ws.Range("J63").Value = 4880
$ws.Range("K63").Value = 10303.3329
$ws.Range("L63").Value = 14640
$ws.Range("M63").Value = -9554.332900000001
$ws.Range("N63").Value = -16138
# row 66
$ws.Range("H66").Value = 3950.7144
$ws.Range("I66").Value = 3434.4443
$ws.Range("J66").Value = 4880
$ws.Range("K66").Value = 30909.9987
$ws.Range("L66").Value = 43920
$ws.Range("M66").Value = -27165.9987
$ws.Range("N66").Value = -51408
# row 113
$ws.Range("H113").Value = 3572026.2
$ws.Range("I113").Value = 609.1667
$ws.Range("K113").Value = 1827.5001
$ws.Range("M113").Value = 342.4999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 62502750
$ws.Range("I80").Value = 62502750
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 62502750
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -62501752
$ws.Range("N80").Value = ""
# row 83
$ws.Range("H83").Value = 62502750
$ws.Range("I83").Value = 62502750
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 312513750
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -312508758
$ws.Range("N83").Value = ""
# row 132
$ws.Range("H132").Value = 5951.3
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 5945.8887
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 17837.6661
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -22897.6661

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 82
$ws.Range("H82").Value = 1251.0344
$ws.Range("I82").Value = 666.5294
$ws.Range("J82").Value = 2079.0833
$ws.Range("K82").Value = 666.5294
$ws.Range("L82").Value = 2079.0833
$ws.Range("M82").Value = -305.5294
$ws.Range("N82").Value = -2801.0833
# row 85
$ws.Range("H85").Value = 1251.0344
$ws.Range("I85").Value = 666.5294
$ws.Range("J85").Value = 2079.0833
$ws.Range("K85").Value = 666.5294
$ws.Range("L85").Value = 2079.0833
$ws.Range("M85").Value = 581.4706
$ws.Range("N85").Value = -4575.0833
# row 132
$ws.Range("H132").Value = 4997.5713
$ws.Range("I132").Value = 3120.5
$ws.Range("J132").Value = 7500.3335
$ws.Range("K132").Value = 9361.5
$ws.Range("L132").Value = 22501.0005
$ws.Range("M132").Value = -6831.5
$ws.Range("N132").Value = -27561.0005
# row 136
$ws.Range("H136").Value = 5326.684
$ws.Range("I136").Value = 1872.4286
$ws.Range("J136").Value = 7341.6665
$ws.Range("K136").Value = 5617.2858
$ws.Range("L136").Value = 22024.9995
$ws.Range("M136").Value = -3067.2858
$ws.Range("N136").Value = -27124.9995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 6025
$ws.Range("I81").Value = 2050
$ws.Range("K81").Value = 4100
$ws.Range("M81").Value = -3039
# row 84
$ws.Range("H84").Value = 6025
$ws.Range("I84").Value = 2050
$ws.Range("K84").Value = 20500
$ws.Range("M84").Value = -15196
# row 136
$ws.Range("H136").Value = 4597.5557
$ws.Range("I136").Value = 3573.4443
$ws.Range("J136").Value = 5109.6113
$ws.Range("K136").Value = 10720.3329
$ws.Range("L136").Value = 15328.8339
$ws.Range("M136").Value = -8170.332900000001
$ws.Range("N136").Value = -20428.8339

